$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(147490.8875300625, 241978.8674148761, 301835.2186333269, 330578.5391073826)
    3  = @(174804.3486792527, 287086.5081288551, 355424.988341677, 392099.6474527332)
    4  = @(150924.3138020896, 262448.1308437714, 334147.6658457189, 370244.3383432403)
    5  = @(114969.0729743874, 190638.6008556033, 234167.2011281125, 259873.451328368)
    6  = @(96060.53792653907, 163154.6752658826, 205673.9356354423, 226086.1117384583)
    7  = @(11262.14893018374, 18119.64078334402, 22249.73253384406, 24086.01157175547)
    8  = @(293532.279030501, 665284.5965777595, 972570.9458413635, 1206414.651504654)
    9  = @(137103.0605244655, 231546.9707504895, 292692.9177647905, 324225.4327157084)
    10 = @(69313.14876375627, 108016.2149676386, 131644.8876273321, 141129.7735666916)
    11 = @(12611.77613632914, 18734.23212412782, 22628.74148001703, 26193.78962240853)
    12 = @(28104.18526028926, 47709.14423609608, 59199.18084798937, 62617.99579466137)
    13 = @(37939.3782314936, 59401.61980814132, 73591.65194541495, 79937.99191583532)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
